$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the old "steps" column (E), shifting it to G
$ws.Columns.Item(5).Insert()
$ws.Columns.Item(5).Insert()

# New header labels for the inserted columns
$ws.Range("E1").Value = "nx,average_clustering(G)"
$ws.Range("F1").Value = "nx.eccentricity(G,starting_node)"

# Copy the header style (bold, border, centered) from D1 onto the new headers
$ws.Range("D1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)

# nx.density(G) values were recomputed; same value for every row
$ws.Range("D2:D51").Value = 0.3057471264367816

# New nx.average_clustering(G) column; same value for every row
$ws.Range("E2:E51").Value = 0.3085360935360936

# New nx.eccentricity(G,starting_node) column (F) and recomputed steps column (G), per row
$ws.Range("F2").Value = 3
$ws.Range("G2").Value = 6
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 13
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 9
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 8
$ws.Range("F6").Value = 2
$ws.Range("G6").Value = 9
$ws.Range("F7").Value = 2
$ws.Range("G7").Value = 9
$ws.Range("F8").Value = 3
$ws.Range("G8").Value = 10
$ws.Range("F9").Value = 3
$ws.Range("G9").Value = 14
$ws.Range("F10").Value = 3
$ws.Range("G10").Value = 9
$ws.Range("F11").Value = 3
$ws.Range("G11").Value = 10
$ws.Range("F12").Value = 3
$ws.Range("G12").Value = 7
$ws.Range("F13").Value = 3
$ws.Range("G13").Value = 8
$ws.Range("F14").Value = 2
$ws.Range("G14").Value = 11
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 9
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 8
$ws.Range("F17").Value = 3
$ws.Range("G17").Value = 13
$ws.Range("F18").Value = 2
$ws.Range("G18").Value = 7
$ws.Range("F19").Value = 3
$ws.Range("G19").Value = 9
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 7
$ws.Range("F21").Value = 3
$ws.Range("G21").Value = 8
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 7
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = 5
$ws.Range("F24").Value = 3
$ws.Range("G24").Value = 8
$ws.Range("F25").Value = 2
$ws.Range("G25").Value = 10
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 9
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 9
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 13
$ws.Range("F29").Value = 3
$ws.Range("G29").Value = 11
$ws.Range("F30").Value = 2
$ws.Range("G30").Value = 10
$ws.Range("F31").Value = 3
$ws.Range("G31").Value = 10
$ws.Range("F32").Value = 3
$ws.Range("G32").Value = 7
$ws.Range("F33").Value = 2
$ws.Range("G33").Value = 8
$ws.Range("F34").Value = 3
$ws.Range("G34").Value = 6
$ws.Range("F35").Value = 3
$ws.Range("G35").Value = 9
$ws.Range("F36").Value = 2
$ws.Range("G36").Value = 8
$ws.Range("F37").Value = 3
$ws.Range("G37").Value = 7
$ws.Range("F38").Value = 3
$ws.Range("G38").Value = 9
$ws.Range("F39").Value = 2
$ws.Range("G39").Value = 5
$ws.Range("F40").Value = 2
$ws.Range("G40").Value = 7
$ws.Range("F41").Value = 3
$ws.Range("G41").Value = 11
$ws.Range("F42").Value = 2
$ws.Range("G42").Value = 13
$ws.Range("F43").Value = 3
$ws.Range("G43").Value = 11
$ws.Range("F44").Value = 2
$ws.Range("G44").Value = 6
$ws.Range("F45").Value = 2
$ws.Range("G45").Value = 8
$ws.Range("F46").Value = 2
$ws.Range("G46").Value = 11
$ws.Range("F47").Value = 3
$ws.Range("G47").Value = 10
$ws.Range("F48").Value = 2
$ws.Range("G48").Value = 14
$ws.Range("F49").Value = 2
$ws.Range("G49").Value = 8
$ws.Range("F50").Value = 3
$ws.Range("G50").Value = 6
$ws.Range("F51").Value = 3
$ws.Range("G51").Value = 11
